# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new column headers, styled like the existing header cells.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting already used by the rest of the header row (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Every player row (2-53) gets the team's season record.
$ws.Range("AD2:AD53").Value = 63
$ws.Range("AE2:AE53").Value = 98
$ws.Range("AF2:AF53").Value = 0
